# Rename the "IdelMine" sheet to "KlingsStaging" and refresh the
# TestCase name cell on Sheet1 that mirrors it, then restore the
# selection recorded on the renamed sheet's tab.

$wb = $excel.ActiveWorkbook

# 1) Rename the second worksheet tab: IdelMine -> KlingsStaging
$ws2 = $wb.Worksheets.Item("IdelMine")
$ws2.Name = "KlingsStaging"

# 2) Sheet1!A2 holds the TestCase name, which tracks the sheet name.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A2").Value = "KlingsStaging"

# 3) Update the selection on the renamed sheet to A4 (it was A2:A40).
$ws2.Activate() | Out-Null
$ws2.Range("A4").Select() | Out-Null
